$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two old "extra-formatted" data rows (old row4 = spa, old row5 = fra). ---
# They carried a wrap-text style (font/fill) that's no longer used after this edit, and a
# custom row height (row4 ht=43.5) that also needs to disappear. Deleting + rewriting clean
# rows below is the simplest way to get fresh, unstyled cells.
$ws.Rows("4:5").Delete()

# --- New header columns (audit columns) appended to row 1 ---
$ws.Range("F1").Value = "cr_by"
$ws.Range("G1").Value = "cr_dtimes"
$ws.Range("H1").Value = "upd_by"
$ws.Range("I1").Value = "upd_dtimes"
$ws.Range("J1").Value = "is_deleted"
$ws.Range("K1").Value = "del_dtimes"

# --- Row 2 (REG / eng) keeps its existing A:E values; just add the audit columns ---
$ws.Range("F2").Value = "System"
$ws.Range("G2").Value = 45526.6013926323
$ws.Range("H2").Value = "NULL"
$ws.Range("I2").Value = "NULL"
$ws.Range("J2").Value = $false
$ws.Range("K2").Value = "NULL"

# --- Row 3 (REG / hin) keeps its existing A:E values; just add the audit columns ---
$ws.Range("F3").Value = "System"
$ws.Range("G3").Value = 45526.6013926323
$ws.Range("H3").Value = "NULL"
$ws.Range("I3").Value = "NULL"
$ws.Range("J3").Value = $false
$ws.Range("K3").Value = "NULL"

# --- Row 4: the former "fra" row (was row 5), rewritten fresh with audit columns ---
$ws.Range("A4").Value = "REG"
$ws.Range("B4").Value = "Régulière"
$ws.Range("C4").Value = "Centre d'inscription régulier"
$ws.Range("D4").Value = "fra"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = "System"
$ws.Range("G4").Value = 45526.6013926323
$ws.Range("H4").Value = "NULL"
$ws.Range("I4").Value = "NULL"
$ws.Range("J4").Value = $false
$ws.Range("K4").Value = "NULL"

# --- Row 5: brand-new test registration-center row ---
$ws.Range("A5").Value = "NEWREG"
$ws.Range("B5").Value = "TEST REGISTRATION CENTER"
$ws.Range("C5").Value = "Testing Purpose"
$ws.Range("D5").Value = "eng"
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = "globaladmin"
$ws.Range("G5").Value = 45636.2292514052
$ws.Range("H5").Value = "vishal"
$ws.Range("I5").Value = 45636.231948358
$ws.Range("J5").Value = $false
$ws.Range("K5").Value = "NULL"

# --- Row 6: the former "spa" row (was row 4), lang_code corrected spa -> es ---
$ws.Range("A6").Value = "REG"
$ws.Range("B6").Value = "Regular"
$ws.Range("C6").Value = "Centro de registro regular"
$ws.Range("D6").Value = "es"
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = "System"
$ws.Range("G6").Value = 45526.6013926323
$ws.Range("H6").Value = "NULL"
$ws.Range("I6").Value = "NULL"
$ws.Range("J6").Value = $false
$ws.Range("K6").Value = "NULL"

# --- Apply the mm:ss.0 (numFmtId 47) date-time number format to the timestamp columns ---
$ws.Range("G2:G6").NumberFormat = "mm:ss.0"
$ws.Range("I5").NumberFormat = "mm:ss.0"

# --- New column D (upd_by) gets an explicit width of 10 ---
$ws.Columns("D").ColumnWidth = 9.1666666667
